$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old sub-header row (Hiver / Eté / Année) - this shifts the
#    10 data rows up by one (old rows 3-12 become rows 2-11).
$ws.Rows.Item(2).Delete()

# 2. Rebuild the header row (row 1).
#    Columns A-E are brand new "index/metadata" columns with no special
#    style (default General style).
$ws.Range("E1").ClearFormats()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Columns F-K keep the existing "Arial 9" font but should end up with a
# cell style that only carries applyFont (no applyNumberFormat) - matching
# fontId=1 numFmtId=0 applyFont="1". Creating a throwaway named cell style
# from font 1 and applying + deleting it leaves exactly that xf behind
# without adding any surviving named style / cellStyleXf.
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$tmpStyle = $wb.Styles.Add("__tmp_header_style__")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "__tmp_header_style__"
$wb.Styles.Item("__tmp_header_style__").Delete()

# 3. Reselect A2:K2, matching the post-edit active selection.
$ws.Range("A2:K2").Select() | Out-Null
